# NV-10 Le Dinh Hau 8-2024.xlsx
# - Fills in the "Đơn sale chính" (personal sale orders) sheet with the
#   August service-order data + a totals row.
# - Propagates the resulting "Chiết khấu sale chính" (main-sale commission)
#   total of 600,000 into the "Lương" (salary) report's SÓC TRĂNG rows and
#   renames the grand-total row to "Tổng lương tại HỆ THỐNG".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$ws2 = $wb.Worksheets.Item("Lương")

# --- Sheet 1: "Đơn sale chính" -------------------------------------------

# Header row
$headers = @(
    "Tiền tố", "Mã dịch vụ", "Ngày thực hiện", "Cơ sở", "Khách hàng",
    "Nguồn khách", "Tên dịch vụ", "Đơn giá gốc", "Sale phụ", "Upsale",
    "Đơn giá", "Đã thanh toán", "Tỉ lệ chiết khấu sale chính",
    "Chiết khấu sale chính"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data row (row 2) - the single service order for the period
$ws1.Cells.Item(2, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(2, 2).Value = 619

# "Ngày thực hiện" must stay literal text ("08-02-2024"), not be
# auto-recognised as a date serial. Force text storage, write the value,
# then drop the cell back to the default "Normal" style so no stray
# number-format is left behind on the cell.
$ws1.Cells.Item(2, 3).NumberFormat = "@"
$ws1.Cells.Item(2, 3).Value = "08-02-2024"
$ws1.Cells.Item(2, 3).Style = "Normal"

$ws1.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$ws1.Cells.Item(2, 5).Value = "mai hồng nương"
$ws1.Cells.Item(2, 6).Value = "Cá nhân"
$ws1.Cells.Item(2, 7).Value = "Thu cánh mũi"
$ws1.Cells.Item(2, 8).Value = 8000000
$ws1.Cells.Item(2, 9).Value = 0
$ws1.Cells.Item(2, 10).Value = 0
$ws1.Cells.Item(2, 11).Value = 8000000
$ws1.Cells.Item(2, 12).Value = 6000000
$ws1.Cells.Item(2, 13).Value = 0.1
$ws1.Cells.Item(2, 14).Value = 600000

# Total row (row 3)
$ws1.Cells.Item(3, 1).Value = "Tổng"
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 3).Value = ""
$ws1.Cells.Item(3, 4).Value = ""
$ws1.Cells.Item(3, 5).Value = ""
$ws1.Cells.Item(3, 6).Value = ""
$ws1.Cells.Item(3, 7).Value = ""
$ws1.Cells.Item(3, 8).Value = 8000000
$ws1.Cells.Item(3, 9).Value = ""
$ws1.Cells.Item(3, 10).Value = 0
$ws1.Cells.Item(3, 11).Value = 8000000
$ws1.Cells.Item(3, 12).Value = 6000000
$ws1.Cells.Item(3, 13).Value = 0
$ws1.Cells.Item(3, 14).Value = 600000

# --- Sheet 2: "Lương" ------------------------------------------------------

# Chiết khấu sale chính tại SÓC TRĂNG
$ws2.Cells.Item(25, 2).Value = 600000
# Tổng lương tại SÓC TRĂNG
$ws2.Cells.Item(34, 2).Value = 600000
# Grand total row: renamed + updated
$ws2.Cells.Item(35, 1).Value = "Tổng lương tại HỆ THỐNG"
$ws2.Cells.Item(35, 2).Value = 600000
